$d = $word.ActiveDocument

# The document ends with a "Requisitos" section listing LOB codes, followed
# by an empty paragraph, a "Ver no Jupiter Salvar em pdf Salvar em docx"
# paragraph and a "(c) 2020 ... Creative Commons Attribution" paragraph
# (site-footer boilerplate), then another empty paragraph and a page-break
# paragraph before the section properties.
#
# The edit removes the footer boilerplate (the empty paragraph right after
# the last LOB requirement line, the "Ver no Jupiter..." paragraph, and the
# "(c) 2020..." paragraph) while leaving the requirement line itself and the
# trailing empty / page-break paragraphs untouched.

$paras = $d.Paragraphs
$reqIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "LOB1215: Recursos*") {
        $reqIndex = $i
    }
}

if ($reqIndex -eq -1) {
    throw "Could not locate the LOB1215 requirement paragraph"
}

# Paragraph right after the requirement line (the blank separator) and the
# paragraph that must be left untouched once the three footer paragraphs
# (blank, "Ver no Jupiter...", "(c) 2020...") are gone.
$deleteStart = $paras.Item($reqIndex + 1).Range.Start
$deleteEnd = $paras.Item($reqIndex + 4).Range.Start

$target = $d.Range($deleteStart, $deleteEnd)

$pattern = [char]13 + "Ver no Jupiter Salvar em pdf Salvar em docx" + [char]13 + `
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution" + [char]13

$replaced = $target.Find.Execute($pattern, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)

if (-not $replaced) {
    throw "Failed to remove the site-footer paragraphs"
}

Write-Output ("Removed footer paragraphs: " + $replaced)
